# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet right after "总计" (pushing 2022-Q2,
# 2022-Q1, 2021-Q4, 2021-Q3, 2020-Q4 each back by one tab position) and
# populates it with the quarter's fund-holdings table, then updates the
# "总计" (totals) sheet with a new leading row for 2022-Q3 (shifting the
# previously-existing rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0, "160644", "鹏华香港美国互联网股票（LOF）人民币", "1.23", "83.13", "3.27", "0.0402", 9),
    @(1, "006792", "鹏华香港美国互联网股票（LOF）美元现汇", "1.23", "83.13", "3.27", "0.0402", 9),
    @(2, "005698", "华夏全球科技先锋混合（QDII）", "0.59", "86.79", "4.21", "0.0248", 9),
    @(3, "015205", "银华全球新能源车量化优选股票（QDII）C", "0.24", "90.37", "3.62", "0.0087", 9),
    @(4, "015204", "银华全球新能源车量化优选股票（QDII）A", "0.19", "90.37", "3.62", "0.0069", 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $idCell = $q3.Cells.Item($rowNum, 1)
    $idCell.Value = $data[0]
    $idCell.Font.Bold = $true
    $idCell.HorizontalAlignment = -4108
    $idCell.VerticalAlignment = -4160
    $idCell.Borders.LineStyle = 1

    # Fund code / name stay text (NumberFormat "@" keeps leading zeros).
    $codeCell = $q3.Cells.Item($rowNum, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $data[1]

    $q3.Cells.Item($rowNum, 3).Value = $data[2]

    for ($c = 3; $c -le 6; $c++) {
        $numCell = $q3.Cells.Item($rowNum, $c + 1)
        $numCell.NumberFormat = "@"
        $numCell.Value = $data[$c]
    }

    $q3.Cells.Item($rowNum, 8).Value = $data[7]
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert the 2022-Q3 row at the top of the
#    data (row 2) and shift the older rows down by one.
# ---------------------------------------------------------------------
$totalRows = @(
    @(0, "2022-Q3", 5, 0.12),
    @(1, "2022-Q2", 6, 0.11),
    @(2, "2022-Q1", 2, 0.19),
    @(3, "2021-Q4", 2, 0.23),
    @(4, "2021-Q3", 2, 0.13),
    @(5, "2020-Q4", 1, 0.01)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $rowNum = $r + 2
    $data = $totalRows[$r]
    $total.Cells.Item($rowNum, 1).Value = $data[0]
    $total.Cells.Item($rowNum, 2).Value = $data[1]
    $total.Cells.Item($rowNum, 3).Value = $data[2]
    $total.Cells.Item($rowNum, 4).Value = $data[3]
}

$total.Range("A1").Select()
